$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BOM")

# --- Update existing rows (prices changed) ---
$ws1.Range("C2").Value = 67
$ws1.Range("C3").Value = 120

# --- Add new BOM rows ---
$ws1.Range("A4").Value = "Бачок"
$ws1.Range("B4").Value = 1
$ws1.Range("C4").Value = 80

$ws1.Range("A5").Value = "Шланг омывателя 5м"
$ws1.Range("B5").Value = 1
$ws1.Range("C5").Value = 50

# Fill the "Sum" formula column, D4:D5 together so it becomes one shared formula
$ws1.Range("D4:D5").Formula = "=B4*C4"

# --- Totals row ---
$ws1.Range("C8").Value = "Итого"
$ws1.Range("D8").Formula = "=SUM(D2:D7)"

# --- Apply the "Calculation" cell style to the whole sum column (D2:D5) ---
$ws1.Range("D2:D5").Style = "Calculation"

# --- Make BOM the active sheet/tab with D9 selected ---
$ws1.Activate() | Out-Null
$ws1.Range("D9").Select() | Out-Null
